$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data (columns F:V) between row 76 and row 77 ---
# Row 76 currently holds the "Zrinjski vs Tuzla City" match; it must become
# the "Siroki Brijeg vs Posusje" match, and row 77 must become the
# "Zrinjski vs Tuzla City" match.

$row76 = @("Siroki Brijeg", 1, "Posusje", 1, 1.9, "04/11/2023 01:13", 2.44, "05/11/2023 12:57", 3.07, "04/11/2023 01:13", 2.69, "05/11/2023 12:57", 3.8, "04/11/2023 01:13", 3.44, "05/11/2023 12:57", "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/siroki-brijeg-posusje/tpgw3CwO/")

$row77 = @("Zrinjski", 3, "Tuzla City", 1, 1.21, "04/11/2023 01:13", 1.16, "05/11/2023 12:44", 5.75, "04/11/2023 01:13", 7.26, "05/11/2023 12:55", 8.65, "04/11/2023 01:13", 14.64, "05/11/2023 12:55", "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/zrinjski-tuzla-city/YJkV4Y8B/")

$cols = @("F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "76").Value = $row76[$i]
    $ws.Range($cols[$i] + "77").Value = $row77[$i]
}

# --- Append a new row 79 with the Posusje vs FK Sarajevo match ---

# Copy formatting from row 78 so the new row matches the existing style
# (bold/bordered index cell, date-formatted match-date cell, ...).
$ws.Range("A78:V78").Copy()
$ws.Range("A79:V79").PasteSpecial(-4122)

$ws.Range("A79").Value = 78
$ws.Range("B79").Value = "bosnia-and-herzegovina"
$ws.Range("C79").Value = "premijer-liga-bih"
$ws.Range("D79").Value = "2023-2024"
$ws.Range("E79").Value = 45241.86458333334
$ws.Range("F79").Value = "Posusje"
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = "FK Sarajevo"
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = 2.67
$ws.Range("K79").Value = "10/11/2023 09:12"
$ws.Range("L79").Value = 2.43
$ws.Range("M79").Value = "11/11/2023 20:42"
$ws.Range("N79").Value = 2.91
$ws.Range("O79").Value = "10/11/2023 09:12"
$ws.Range("P79").Value = 3.14
$ws.Range("Q79").Value = "11/11/2023 20:33"
$ws.Range("R79").Value = 2.58
$ws.Range("S79").Value = "10/11/2023 09:12"
$ws.Range("T79").Value = 2.53
$ws.Range("U79").Value = "11/11/2023 20:42"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/posusje-fk-sarajevo/j9uaKD8b/"

$wb.Save()
